{"js": "// Update the stack trace text in the body to match the new test harness\n// line numbers / generated-accessor index, and swap the Maven\n// Surefire/Tycho/Equinox launch frames for the Eclipse JDT JUnit runner\n// frames (test template version bump).\n\nconst body = context.document.body;\n\nconst oldBlock1 =\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)\\n\" +\n  \"\\tat sun.reflect.GeneratedMethodAccessor5.invoke(Unknown Source)\";\nconst newBlock1 =\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)\\n\" +\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)\\n\" +\n  \"\\tat sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)\";\n\nconst oldBlock2 =\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)\\n\" +\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)\\n\" +\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)\\n\" +\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)\\n\" +\n  \"\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:161)\\n\" +\n  \"\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)\\n\" +\n  \"\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)\\n\" +\n  \"\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)\\n\" +\n  \"\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)\\n\" +\n  \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)\\n\" +\n  \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)\\n\" +\n  \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)\\n\" +\n  \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)\";\nconst newBlock2 =\n  \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\";\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for block, got \" + results.items.length\n    );\n  }\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(oldBlock1, newBlock1);\nawait replaceOnce(oldBlock2, newBlock2);\n", "ps1": "# Update the stack trace text in the body to match the new test harness\n# line numbers / generated-accessor index, and swap the Maven\n# Surefire/Tycho/Equinox launch frames for the Eclipse JDT JUnit runner\n# frames (test template version bump).\n\n$d = $word.ActiveDocument\n\n$oldBlock1 = \"`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)`n\" + `\n             \"`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)`n\" + `\n             \"`tat sun.reflect.GeneratedMethodAccessor5.invoke(Unknown Source)\"\n$newBlock1 = \"`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)`n\" + `\n             \"`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)`n\" + `\n             \"`tat sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)\"\n\n$oldBlock2 = \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n\" + `\n             \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n\" + `\n             \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n\" + `\n             \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n\" + `\n             \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n             \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n             \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n             \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n             \"`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n\" + `\n             \"`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n\" + `\n             \"`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n\" + `\n             \"`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n\" + `\n             \"`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n\" + `\n             \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n             \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n             \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n             \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n             \"`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n\" + `\n             \"`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n\" + `\n             \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n\" + `\n             \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n\" + `\n             \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n\" + `\n             \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n\" + `\n             \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n             \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n             \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n             \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n             \"`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n\" + `\n             \"`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n\" + `\n             \"`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n\" + `\n             \"`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)\"\n$newBlock2 = \"`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n\" + `\n             \"`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n\" + `\n             \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n\" + `\n             \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n\" + `\n             \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n\" + `\n             \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldBlock1\n$find.Replacement.Text = $newBlock1\n$find.MatchCase = $true\n$find.Execute($oldBlock1, $true, $false, $false, $false, $false, $true, 1, $false, $newBlock1, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = $oldBlock2\n$find2.Replacement.Text = $newBlock2\n$find2.MatchCase = $true\n$find2.Execute($oldBlock2, $true, $false, $false, $false, $false, $true, 1, $false, $newBlock2, 2) | Out-Null\n"}
